$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source data shrank from 4 fines to 2 fines - drop the last two data rows.
$ws.Rows.Item(5).Delete()
$ws.Rows.Item(4).Delete()

# Plate Number (col C) and Fine Number (col H) are numeric-looking text values
# that must stay stored as text, not numbers. Force Text format only for the
# duration of the write, then restore the Normal style so no stray formatting
# is left behind on the cells.
$ws.Range("C2:C3").NumberFormat = "@"
$ws.Range("H2:H3").NumberFormat = "@"

# Row 2: previously the FORD MUSTANG fine -> now the CADILLAC ESCALADE fine
$ws.Range("A2").Value = "CADILLAC ESCALADE, 2023, Blue"
$ws.Range("B2").Value = "N"
$ws.Range("C2").Value = "85540"
$ws.Range("D2").Value = "14 Jul 2025, 12:10 am"
$ws.Range("E2").Value = "Ras Al khour St"
$ws.Range("F2").Value = "Dubai Police"
$ws.Range("G2").Value = "AED 600"
$ws.Range("H2").Value = "7037866556"
$ws.Range("I2").Value = "Exceeding maximum speed limit by not more than 30 km h"
$ws.Range("J2").Value = "Please contact Dubai Police for details about disputing your fine."

# Row 3: the KIA K5 fine, with its Amount/Details normalized to match row 2
$ws.Range("A3").Value = "KIA K5, 2023, Black"
$ws.Range("B3").Value = "DD"
$ws.Range("C3").Value = "81392"
$ws.Range("D3").Value = "11 Jul 2025, 8:30 am"
$ws.Range("E3").Value = "Dubai Alain Road"
$ws.Range("F3").Value = "Dubai Police"
$ws.Range("G3").Value = "AED 600"
$ws.Range("H3").Value = "7037841032"
$ws.Range("I3").Value = "Exceeding maximum speed limit by not more than 30 km h"
$ws.Range("J3").Value = "Please contact Dubai Police for details about disputing your fine."

# Restore the default (Normal) style on the reformatted cells so no text
# number-format lingers on them.
$ws.Range("C2:C3").Style = "Normal"
$ws.Range("H2:H3").Style = "Normal"
